# Update the division problems in the single table of the worksheet.
# Each cell's text is replaced in place (row/col addressed) so existing
# run formatting (TimeNewRoman, sz 30, left-justified paragraph) is kept
# and no Find/Replace ambiguity arises from values that coincide between
# the "before" and "after" sets (e.g. "37÷6=").

$d = $word.ActiveDocument
$t = $d.Tables(1)

$updates = @(
    @{ Row = 1;  Col = 1; Old = "52÷7="; New = "66÷4=" },
    @{ Row = 1;  Col = 2; Old = "89÷7="; New = "55÷3=" },
    @{ Row = 1;  Col = 3; Old = "88÷8="; New = "85÷4=" },
    @{ Row = 1;  Col = 4; Old = "31÷3="; New = "41÷5=" },
    @{ Row = 1;  Col = 5; Old = "56÷4="; New = "51÷7=" },

    @{ Row = 5;  Col = 1; Old = "23÷3="; New = "39÷6=" },
    @{ Row = 5;  Col = 2; Old = "62÷2="; New = "89÷2=" },
    @{ Row = 5;  Col = 3; Old = "77÷7="; New = "94÷2=" },
    @{ Row = 5;  Col = 4; Old = "52÷2="; New = "12÷2=" },
    @{ Row = 5;  Col = 5; Old = "49÷2="; New = "34÷6=" },

    @{ Row = 9;  Col = 1; Old = "25÷2="; New = "32÷4=" },
    @{ Row = 9;  Col = 2; Old = "37÷6="; New = "21÷6=" },
    @{ Row = 9;  Col = 3; Old = "37÷4="; New = "37÷6=" },
    @{ Row = 9;  Col = 4; Old = "22÷8="; New = "13÷7=" },
    @{ Row = 9;  Col = 5; Old = "97÷2="; New = "81÷6=" },

    @{ Row = 13; Col = 1; Old = "57÷8="; New = "81÷4=" },
    @{ Row = 13; Col = 2; Old = "72÷7="; New = "65÷6=" },
    @{ Row = 13; Col = 3; Old = "36÷8="; New = "74÷3=" },
    @{ Row = 13; Col = 4; Old = "35÷7="; New = "60÷9=" },
    @{ Row = 13; Col = 5; Old = "14÷8="; New = "78÷7=" },

    @{ Row = 17; Col = 1; Old = "73÷9="; New = "29÷3=" },
    @{ Row = 17; Col = 2; Old = "36÷3="; New = "78÷7=" },
    @{ Row = 17; Col = 3; Old = "92÷9="; New = "60÷7=" },
    @{ Row = 17; Col = 4; Old = "41÷7="; New = "84÷2=" },
    @{ Row = 17; Col = 5; Old = "92÷3="; New = "24÷7=" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $range = $cell.Range
    # Trim the trailing cell-mark (CR + cell-end char) that Range.Text
    # includes so we only compare/replace the visible content.
    $current = $range.Text.Substring(0, $range.Text.Length - 2)
    if ($current -ne $u.Old) {
        Write-Host "WARNING: cell ($($u.Row),$($u.Col)) expected '$($u.Old)' but found '$current'"
    }
    $range.Text = $u.New
}

Write-Host "Done updating $($updates.Count) cells."
